# Actualización desde MV -datos-
# Appends 5 new daily rows (04-10-2021 .. 08-10-2021) to the bottom of the
# "Swap promedio camara" sheet, matching the source data feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Date = "04-10-2021"; Values = @(4.53, 4.8,  5.01, 5.19, 5.69, -0.86, 0.48, 0.95, 1.34, 1.68, 2.31, 2.48) },
    @{ Date = "05-10-2021"; Values = @(4.66, 4.96, 5.17, 5.37, 5.91, -0.75, 0.59, 1.08, 1.49, 1.84, 2.54, 2.72) },
    @{ Date = "06-10-2021"; Values = @(4.68, 5,    5.21, 5.43, 5.98, -0.81, 0.61, 1.13, 1.53, 1.89, 2.6,  2.77) },
    @{ Date = "07-10-2021"; Values = @(4.69, 5.01, 5.23, 5.43, 5.97, -0.89, 0.55, 1.1,  1.51, 1.87, 2.59, 2.76) },
    @{ Date = "08-10-2021"; Values = @(4.99, 5.3,  5.49, 5.67, 6.17, -1.02, 0.54, 1.12, 1.55, 1.95, 2.74, 2.94) }
)

$startRow = 194

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $entry = $newRows[$i]

    # Column A holds a date-like label ("04-10-2021") that must stay plain
    # text (matching the rest of the column) instead of being auto-parsed
    # into a date serial by Excel's input heuristics. Writing it as a text
    # formula and then converting the formula to its value in place keeps
    # it a genuine shared-string text cell with no extra cell formatting.
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Formula = '="' + $entry.Date + '"'
    $dateCell.Copy()
    $dateCell.PasteSpecial(-4163) # xlPasteValues

    $values = $entry.Values
    for ($c = 0; $c -lt $values.Count; $c++) {
        $ws.Cells.Item($r, 2 + $c).Value = $values[$c]
    }
}

$excel.CutCopyMode = 0
